# Update gh-pages to output generated at 456a3b4
# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 12642
    3  = 604
    5  = 16
    6  = 277
    7  = 393
    9  = 12617
    11 = 3128
    12 = 541
    19 = 657
    20 = 2840
    21 = 6112
    23 = 3611
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates2 = @{
    2  = 12642
    3  = 604
    5  = 16
    6  = 277
    8  = 393
    10 = 12617
    12 = 3128
    13 = 541
    20 = 657
    21 = 2840
    23 = 6112
    25 = 3611
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates2.Keys) {
    $ws4.Range("F$row").Value = $updates2[$row]
}
